$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of data (row 15): ITEM 13, CDI SEMILLITAS, 24/02/2019 (serial 43520),
#     VALLE DEL CAUCA, YOTOCO ---------------------------------------------------
# Write the new shared-string-bearing cells in the same order the strings were
# first introduced (VALLE DEL CAUCA, YOTOCO, CDI SEMILLITAS) so the shared
# string table ends up in the expected order.
$ws.Range("D15").Value = "VALLE DEL CAUCA "
$ws.Range("E15").Value = "YOTOCO"
$ws.Range("B15").Value = "CDI SEMILLITAS"
$ws.Range("A15").Value = 13
$ws.Range("C15").Value = 43520

# Match the formatting already used by the other filled-in rows (e.g. row 12):
# thin border on every side, 11pt Calibri, no fill.
$ws.Range("B15").Borders.LineStyle = 1
$ws.Range("B15").Font.Size = 11
$ws.Range("D15").Borders.LineStyle = 1
$ws.Range("D15").Font.Size = 11
$ws.Range("E15").Borders.LineStyle = 1
$ws.Range("E15").Font.Size = 11

# --- Column widths for I, J, K ------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 16.166666666666668
$ws.Columns.Item(10).ColumnWidth = 23.5
$ws.Columns.Item(11).ColumnWidth = 22

# --- Selection moves to the new row (entire row 9 selected) ------------------
$ws.Rows.Item(9).Select() | Out-Null

# --- Page setup: landscape, no explicit paper size ----------------------------
$ws.PageSetup.PaperSize = $null
$ws.PageSetup.Orientation = 2
